$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in row 2
$ws.Range("B2").Value = 5

# Update values in row 3
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 1

# Update values in row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 1

# Remove row 5 entirely (was A5=0, B5=1)
$ws.Rows.Item(5).Delete()
